$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the missing "Horas reales" value for the existing last task row (49)
$ws.Range("C49").Value = 2

# Add a new task row 50: "Max gaps por dia"
$ws.Range("A50").Value = "Max gaps por dia"
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = "Clara"
$ws.Range("E50").Value = 42701

# Match the date formatting used in column E for the other task rows
$ws.Range("E49").Copy()
$ws.Range("E50").PasteSpecial(-4122)
$ws.Range("E50").Value = 42701
$excel.CutCopyMode = $false

# Update the view to reflect the scrolled/selected state from the edit
$ws.Activate()
$ws.Range("F51").Select()
